$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-51: Row, Coin (B), Link (C), Price (D), Volume1h (E)
$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '26.036.27', '  +0.82%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.747.38', '  +0.38%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9991', '  -0.12%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '234.19', '  +3.31%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9995', '  -0.05%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.5293', '  +2.41%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2784', '  +2.11%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06184', '  +1.50%  '),
    @(10, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.744.79', '  +0.22%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07239', '  +3.45%  '),
    @(12, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '15.34', '  +1.12%  '),
    @(13, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6436', '  +2.03%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.603', '  +2.51%  '),
    @(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '78.34', '  +2.60%  '),
    @(16, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9992', '  -0.09%  '),
    @(17, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9988', '  -0.11%  '),
    @(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '25.957.45', '  +0.42%  '),
    @(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '11.62', '  +1.52%  '),
    @(20, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000006742', '  +1.88%  '),
    @(21, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.969.88', '  +0.36%  '),
    @(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '4.315', '  +6.10%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.816', '  +4.88%  '),
    @(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '5.210', '  +2.43%  '),
    @(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '139.30', '  +1.90%  '),
    @(26, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.513', '  +0.72%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '15.31', '  +2.33%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.809', '  -0.41%  '),
    @(29, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '104.19', '  +1.51%  '),
    @(30, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.08321', '  +0.19%  '),
    @(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '3.796', '  +4.89%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '3.674', '  +8.72%  '),
    @(33, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04540', '  +3.22%  '),
    @(34, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.643', '  +0.96%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.000', '  +3.23%  '),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6324', '  +6.16%  '),
    @(37, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.709', '  +1.20%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01596', '  +2.78%  '),
    @(39, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.938', '  +0.23%  '),
    @(40, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9990', '  -0.05%  '),
    @(41, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '98.55', '  -2.63%  '),
    @(42, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3909', '  +2.41%  '),
    @(43, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.7396', '  +2.07%  '),
    @(44, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '5.042', '  +3.56%  '),
    @(45, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1143', '  +3.93%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '6.337', '  +2.49%  '),
    @(47, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05346', '  -2.49%  '),
    @(48, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '54.00', '  +4.00%  '),
    @(49, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '30.70', '  +3.28%  '),
    @(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '7.695', '  +3.57%  '),
    @(51, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.3466', '  +2.32%  ')
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]

    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item[3]
    $dCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $item[4]
}
